$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Default name (row 2) changes from FlameSpirit to FlameSpirit_Easy
$ws.Range("B2").Value = "FlameSpirit_Easy"

# Row 10: Easy difficulty
$ws.Range("B10").Value = "FlameSpirit_Easy"
$ws.Range("E10").Value = 300
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 6

# Row 11: Normal difficulty (new data added)
$ws.Range("B11").Value = "FlameSpirit_Normal"
$ws.Range("E11").Value = 250
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 5

# Row 12: Hard difficulty (new row of data)
$ws.Range("B12").Value = "FlameSpirit_Hard"
$ws.Range("E12").Value = 200
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 3

# Update the active selection to H11 to match the saved view state
$ws.Range("H11").Select()
